$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 06:22"

# Row 8 - Alemania
$ws.Range("D8").Value = 99400
$ws.Range("E8").Value = 43967

# Row 14 - Brasil
$ws.Range("B14").Value = 43368
$ws.Range("C14").Value = 289
$ws.Range("E14").Value = 16282
$ws.Range("G14").Value = 20
$ws.Range("H14").Value = 2761

# Row 43 - Australia
$ws.Range("D43").Value = 4912
$ws.Range("E43").Value = 1661

# Row 64 - Kazajistan
$ws.Range("B64").Value = 2025
$ws.Range("C64").Value = 30
$ws.Range("E64").Value = 1517

# Row 117 - Sri Lanka
$ws.Range("F117").Value = 2

# Row 127 - Paraguay
$ws.Range("E127").Value = 142
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 9

# Row 130 - Martinica
$ws.Range("B130").Value = 164
$ws.Range("C130").Value = 1
$ws.Range("E130").Value = 77
$ws.Range("F130").Value = 6

# Row 145 - Guayana Francesa
$ws.Range("D145").Value = 83
$ws.Range("E145").Value = 13

# Row 166 - Nepal
$ws.Range("D166").Value = 5
$ws.Range("E166").Value = 37

# Row 172 - Mongolia
$ws.Range("B172").Value = 35
$ws.Range("C172").Value = 1
$ws.Range("E172").Value = 27

# Row 188 - Dominica
$ws.Range("D188").Value = 9
$ws.Range("E188").Value = 7

$wb.Save()
